# Update the cryptos list worksheet with the latest scraped values
# (GitHub Actions refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17/18, 32/33 and 38/39 swapped their coin ordering (ranking
# changed), so Coin (B) / Link (C) are rewritten too for those rows.
# All other rows only refresh Price (D) and/or Volume(1h) (E).
#
# Price values that look like plain decimals (e.g. "604.73", "0.750")
# would otherwise be auto-coerced to numbers by Excel's normal typing
# behaviour, which both changes their cell type and mangles the value
# (trailing zeros dropped, binary-float rounding noise). To keep them
# as the literal text the source data has, they're entered with a
# leading apostrophe (Excel's standard "treat as text" prefix) and the
# cell style is then reset to Normal so no visible/structural style
# change remains - only the underlying value type becomes Text,
# matching the original inline-string cells.

function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Price (D) values that could otherwise be misread as numbers.
$textPrices = @{
    "D5" = "604.73"
    "D6" = "157.81"
    "D9" = "0.549"
    "D12" = "0.504"
    "D14" = "39.31"
    "D17" = "7.33"
    "D20" = "508.59"
    "D21" = "15.40"
    "D22" = "0.750"
    "D23" = "8.09"
    "D24" = "14.71"
    "D25" = "86.25"
    "D26" = "0.172"
    "D29" = "9.11"
    "D30" = "2.36"
    "D31" = "6.91"
    "D32" = "28.30"
    "D33" = "2.85"
    "D36" = "6.39"
    "D38" = "3.34"
    "D39" = "55.34"
    "D40" = "492.54"
    "D41" = "0.0428"
    "D42" = "0.129"
    "D43" = "8.78"
    "D45" = "2.48"
    "D47" = "28.41"
    "D48" = "2.42"
}

foreach ($addr in $textPrices.Keys) {
    Set-TextValue $addr $textPrices[$addr]
}

$rows = @{
    2 = @{ D = "66.725.42"; E = "  +0.50%  " }
    3 = @{ D = "3.245.05"; E = "  +1.80%  " }
    4 = @{ E = "  -0.09%  " }
    5 = @{ E = "  +0.36%  " }
    6 = @{ E = "  +1.54%  " }
    7 = @{ E = "  -0.01%  " }
    8 = @{ D = "3.244.70"; E = "  +1.73%  " }
    9 = @{ E = "  +0.06%  " }
    10 = @{ E = "  +2.25%  " }
    11 = @{ E = "  +0.17%  " }
    12 = @{ E = "  -0.89%  " }
    13 = @{ E = "  +2.96%  " }
    14 = @{ E = "  +1.26%  " }
    15 = @{ D = "3.773.69"; E = "  +1.60%  " }
    16 = @{ D = "66.709.25"; E = "  +0.36%  " }
    17 = @{ B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; E = "  -0.90%  " }
    18 = @{ B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.242.76"; E = "  +1.67%  " }
    19 = @{ E = "  +1.69%  " }
    20 = @{ E = "  -0.68%  " }
    21 = @{ E = "  -0.40%  " }
    22 = @{ E = "  +2.39%  " }
    23 = @{ E = "  -0.42%  " }
    24 = @{ E = "  -1.13%  " }
    25 = @{ E = "  +1.98%  " }
    26 = @{ E = "  +92.41%  " }
    27 = @{ E = "  +0.13%  " }
    28 = @{ E = "  +0.42%  " }
    29 = @{ E = "  -0.72%  " }
    30 = @{ E = "  -1.22%  " }
    31 = @{ E = "  -2.14%  " }
    32 = @{ B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; E = "  +0.71%  " }
    33 = @{ B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; E = "  -6.86%  " }
    34 = @{ E = "  -0.02%  " }
    35 = @{ E = "  -4.39%  " }
    36 = @{ E = "  -2.42%  " }
    37 = @{ D = "0.0₃0801"; E = "  +19.00%  " }
    38 = @{ B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; E = "  +18.11%  " }
    39 = @{ B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; E = "  +1.08%  " }
    40 = @{ E = "  -3.50%  " }
    41 = @{ E = "  +1.22%  " }
    42 = @{ E = "  +2.19%  " }
    43 = @{ E = "  -0.95%  " }
    44 = @{ E = "  -2.01%  " }
    45 = @{ E = "  +1.71%  " }
    46 = @{ D = "2.944.84"; E = "  +3.35%  " }
    47 = @{ E = "  +0.13%  " }
    48 = @{ E = "  -0.08%  " }
    49 = @{ E = "  +2.01%  " }
    50 = @{ E = "  +0.01%  " }
    51 = @{ E = "  +0.00%  " }
}

foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
